$wb = $excel.ActiveWorkbook

$wsDeals   = $wb.Worksheets.Item(1)   # Deals_data
$wsBanking = $wb.Worksheets.Item(2)   # Banking_Details
$wsSheet1  = $wb.Worksheets.Item(3)   # Sheet1

# ------------------------------------------------------------------
# 1. Update values on the "Deals_data" sheet (rows 3 & 4) so the test
#    deal references a new FMT Local Test Deal / Test Bundle 3
#    Provider / new MSISDN / new city, and apply the "code" style
#    (Courier New, already used for the equivalent rows on Sheet1)
#    to the providerName/dealName columns for rows 2-4.
# ------------------------------------------------------------------

# Copy the existing "code" formatting (style used in Sheet1 rows 8-10)
# onto Deals_data B2:C4 without touching the values.
$wsSheet1.Range("B8:C10").Copy()
$wsDeals.Range("B2:C4").PasteSpecial(-4122)

# Now update the actual cell contents that changed.
$wsDeals.Range("C3").Value = "FMT Local Test Deal"
$wsDeals.Range("B4").Value = "Test Bundle 3 Provider"
$wsDeals.Range("C4").Value = "FmtLocalEssentialsDeal_Safebase1_Bundle_DealDescription"
$wsDeals.Range("D4").Value = "'0678678771"
$wsDeals.Range("E4").Value = "Cape Town"

# ------------------------------------------------------------------
# 2. Update the saved selection / active sheet state.
#    Final state: Deals_data is the active tab with B2 selected;
#    Banking_Details keeps D2 selected; Sheet1 keeps A8:E10 selected.
#    Select the non-active sheets first so the very last Select()
#    call (on Deals_data) is the one that ends up "active".
# ------------------------------------------------------------------

$wsBanking.Range("D2").Select()
$wsSheet1.Range("A8:E10").Select()
$wsDeals.Range("B2").Select()
$wsDeals.Activate()
